# Apply "Tested all photos on log" edit to eyeDetectionTestingLog workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 0. Stash pristine copies of the two "A column" look-and-feels we will
#    need later (plain wrap/shrink style used by rows without a
#    hyperlink-looking font, and the Hyperlink-font wrap/shrink style)
#    in scratch cells far outside the used range. Re-applying a
#    Hyperlinks.Add() always perturbs a cell's style slightly, so we
#    keep a clean source to re-paste from afterwards.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 1. Remove all existing hyperlinks up-front: hyperlink Range objects in
#    this engine do not shift when rows are deleted, so we rebuild the
#    whole collection from scratch once the sheet has its final shape.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 2. Remove the duplicate row (old row 10 repeated the text of row 9),
#    shifting every following row up by one.
# ---------------------------------------------------------------------
$ws.Rows(10).Delete()

# ---------------------------------------------------------------------
# 3. Give every data row (4-19) the same B/C/D/E number formatting that
#    rows 2-3 already use (text wrap for B/C/E, short date for D) before
#    filling in the actual values.
# ---------------------------------------------------------------------
$ws.Range("B2:E2").Copy()
$ws.Range("B4:E19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. Fill in the testing-log results for every photo row.
# ---------------------------------------------------------------------
$rows = @{
    4  = @("Success", "None")
    5  = @("Error",   "Nostril and right ear falsely recognized")
    6  = @("Error",   "Some false recognition on left forehead. Picture is too large to see the rest of and cannot be scrolled across or resized")
    7  = @("?",        "Cannot see eyes because window can't be resized")
    8  = @("Error",   "Left eye double selected")
    9  = @("Success", "None")
    10 = @("Success", "None")
    11 = @("Error",   "Doesn't even recognize the face correctly")
    12 = @("Error",   "Falsely recognizes left nostril")
    13 = @("Error",   "Does not recognize eyes")
    14 = @("Error",   "Double selects left eye, does not select right eye")
    15 = @("Success", "None")
    16 = @("Error",   "Traceback (most recent call last):`n  File `"C:\Users\Shannon\Documents\GitHub\DVS-Python\eyeDetection.py`", line 99, in <module>`n    image = DetectRedEyes(img, faceCascade, eyeCascade)`n  File `"C:\Users\Shannon\Documents\GitHub\DVS-Python\eyeDetection.py`", line 71, in DetectRedEyes`n    cv.SetImageROI(image, (pt1[0],`nUnboundLocalError: local variable 'pt1' referenced before assignment")
    17 = @("Error",   "Does not recognize eyes")
    18 = @("Error",   "Does not capture full right eye")
    19 = @("Success", "None")
}

foreach ($r in 4..19) {
    $vals = $rows[$r]
    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = 41389
    $ws.Range("E$r").Value = "Shannon Harris"
}

# ---------------------------------------------------------------------
# 5. Row 16 (old row 17, "africa_faces" traceback entry) needs a taller
#    row to fit the full traceback text; every other row kept its
#    original auto-fit height when rows shifted up.
# ---------------------------------------------------------------------
$ws.Rows(16).RowHeight = 210

# ---------------------------------------------------------------------
# 6. Re-create the hyperlinks. The ones that existed before keep
#    pointing at the same web addresses (now shifted up one row because
#    of the deleted duplicate); six rows that never had a clickable URL
#    before now also get one.
# ---------------------------------------------------------------------
$links = @{
    3  = "https://onlinehealthsafe.com/ohs/wp-content/uploads/2012/12/face.jpg"
    4  = "http://upload.wikimedia.org/wikipedia/commons/thumb/5/55/Mona_Lisa_headcrop.jpg/250px-Mona_Lisa_headcrop.jpg"
    5  = "http://media.npr.org/assets/img/2012/05/30/rudyeugene_custom-b19998a160cf11e6274e774f6228ed385d8c7126-s6-c10.jpg"
    6  = "http://upload.wikimedia.org/wikipedia/commons/e/e7/Boy_Face_from_Venezuela.jpg"
    7  = "http://www.hypergridbusiness.com/wp-content/uploads/2010/05/dshiao_headshot2.jpg"
    8  = "http://www.huhmagazine.co.uk/images/uploaded/pittchanel_01.jpg"
    9  = "http://face2face.si.edu/.a/6a00e550199efb8833010536a5483e970c-800wi"
    16 = "http://news.nationalgeographic.com/news/2004/10/photogalleries/africa_faces/images/primary/faces_p9.jpg"
    18 = "http://3.bp.blogspot.com/_E_Emkyb959E/TBqJ6IAOOBI/AAAAAAAAAOc/06Rx3egsgPU/s1600/Face.JPG"
    19 = "http://farm3.staticflickr.com/2060/2518721774_6634aaaa02.jpg"
    2  = "http://bloximages.chicago2.vip.townnews.com/azstarnet.com/content/tncms/assets/v3/editorial/3/78/3781cbf8-0d81-511b-8edd-10fcbf43ace5/50ef7ffaab184.preview-620.jpg"
    10 = "http://0.tqn.com/d/menshair/1/0/k/5/-/-/round-black.jpg"
    12 = "http://4.bp.blogspot.com/_c7BpG_CfMJE/S7GhAKIMujI/AAAAAAAAABE/hzGx1O6ZSkI/s1600/smiling-faces01.jpg"
    13 = "http://img.izismile.com/img/img3/20100901/640/wrinkled_faces_640_12.jpg"
    14 = "http://corkap.files.wordpress.com/2011/03/faces022b.jpg"
    15 = "http://www.frugaltravelguy.com/wp-content/uploads/2012/05/Faces-of-FlyerTalk-cdking-Headshot-Walking-by-the-Charles-River-Boston.jpg"
    17 = "http://www.simonhoegsberg.com/faces_of_new_york/images/01_faces.jpg"
}
$linkOrder = @(3, 4, 5, 6, 7, 8, 9, 16, 18, 19, 2, 10, 12, 13, 14, 15, 17)

foreach ($r in $linkOrder) {
    $ws.Hyperlinks.Add($ws.Range("A$r"), $links[$r])
}

# ---------------------------------------------------------------------
# 7. Adding a hyperlink always nudges the cell's style, so restore the
#    exact look every A-column cell should have: the plain hyperlink
#    wrap/shrink style for rows that used it before (2,4,5,6,7,8,9,16,
#    18,19) and the same style for the six rows that are newly linked
#    (10,12,13,14,15,17) too - row 3 already used it unmodified.
# ---------------------------------------------------------------------
$ws.Range("H1").Copy()
foreach ($r in @(2, 4, 5, 6, 7, 8, 9, 16, 18, 19)) {
    $ws.Range("A$r").PasteSpecial(-4122)
}
$ws.Range("H2").Copy()
foreach ($r in @(3, 10, 12, 13, 14, 15, 17)) {
    $ws.Range("A$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Remove the scratch cells so they don't show up in the saved sheet.
$ws.Range("H1:H2").Clear()

# ---------------------------------------------------------------------
# 8. Sheet view: selection now sits on B19, and the window is scrolled
#    back to the top-left corner (no more frozen A3 top-left cell).
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B19").Select()
